$d = $word.ActiveDocument

# --- Change 1 --------------------------------------------------------
# The "Igår så började vi programmera..." paragraph had two words
# ("griden" / "gridden") wrapped in spellcheck proofErr runs, splitting
# the paragraph into several <w:r> elements. Re-typing the identical
# text over the whole paragraph collapses it back into a single run
# and drops the proofErr markers.
$old1 = "Igår så började vi programmera i java Dennis började med att skapa fönstret och den generella griden, och Emma började med registreringen. Idag ska Emma se till så att formuläret fungerar. Dennis ska fina till den generella gridden och möjligt vi börja med att skapa händelser."
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $old1, 2) | Out-Null

# --- Change 2 --------------------------------------------------------
# Append three new paragraphs (2016-11-29 entry) right after the
# "Vi får fortsätta..." paragraph. A trailing sentinel character is
# appended after the new text so the document's final paragraph-mark
# boundary is not the exact insertion point (that edge position trips
# up Bookmarks.Add); the sentinel is stripped immediately afterwards.
$old2 = "Vi får fortsätta med det vi skulle gjort i lördags. Pontus ska göra knappar till vyn, men det ska inte fungera än."
$new2 = "Vi får fortsätta med det vi skulle gjort i lördags. Pontus ska göra knappar till vyn, men det ska inte fungera än.^p2016-11-29^pIgår höll Dennis på med händelsehanterare, emma fortsatte med registreringen och Pontus började med navigering.^pIdag får vi se hur mycket vi hinner. 9.45 ska vi gå till Combitech på föreläsning om säkerhet, och sedan ska vi fortsätta om tid finns. Om inte annat så får vi fortsätta imorgon.#"
$d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# The "_GoBack" bookmark used to sit at the end of the (now earlier)
# "Vi får fortsätta..." paragraph; move it so it again marks the very
# end of the document content, after the newly appended text.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$endPos = $d.Content.End
$newBmPos = $endPos - 2
$r = $d.Range($newBmPos, $newBmPos)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null

# Remove the temporary sentinel character now that the bookmark is
# anchored just before it (it survives the delete and lands right
# after the real text, before the paragraph mark).
$sentinelRange = $d.Range($endPos - 2, $endPos - 1)
$sentinelRange.Delete()
